# Update cryptocurrency price/volume data in Sheet1 to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.781.83'
$ws.Range("E2").Value = '  +6.43%  '
$ws.Range("D3").Value = '1.737.32'
$ws.Range("E3").Value = '  +5.10%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.44'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5467'
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2764'
$ws.Range("E8").Value = '  +3.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06731'
$ws.Range("E9").Value = '  +5.62%  '
$ws.Range("E10").Value = '  +6.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07784'
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.687'
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.975.95'
$ws.Range("E13").Value = '  +5.10%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.723.93'
$ws.Range("E14").Value = '  +5.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5991'
$ws.Range("E15").Value = '  +6.46%  '
$ws.Range("D16").Value = '0.0₅8434'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.43'
$ws.Range("E17").Value = '  +5.92%  '
$ws.Range("D18").Value = '27.788.86'
$ws.Range("E18").Value = '  +6.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.63'
$ws.Range("E19").Value = '  +18.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.835'
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  +5.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.231'
$ws.Range("E23").Value = '  +4.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.66'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1250'
$ws.Range("E26").Value = '  +3.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.711'
$ws.Range("E27").Value = '  +13.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.456'
$ws.Range("E28").Value = '  +2.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.17'
$ws.Range("E29").Value = '  +7.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05660'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.312'
$ws.Range("E31").Value = '  +3.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.693'
$ws.Range("E32").Value = '  +5.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.516'
$ws.Range("E33").Value = '  +3.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.679'
$ws.Range("E34").Value = '  +6.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9767'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.856'
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5957'
$ws.Range("E38").Value = '  +3.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01670'
$ws.Range("E39").Value = '  +4.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.882'
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8491'
$ws.Range("E41").Value = '  +0.59%  '
$ws.Range("D42").Value = '1.049.56'
$ws.Range("E42").Value = '  +2.79%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.16'
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = '1.881.13'
$ws.Range("E45").Value = '  +5.04%  '
$ws.Range("E46").Value = '  +11.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.45'
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.257'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("E49").Value = '  +2.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9996'
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05313'
$ws.Range("E51").Value = '  -0.53%  '
